$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "684×8=5472" "423×7=2961"
Replace-Text "149×2=298" "101×6=606"
Replace-Text "848×5=4240" "113×2=226"
Replace-Text "824×8=6592" "312×5=1560"
Replace-Text "217×7=1519" "720×7=5040"
Replace-Text "488×7=3416" "839×8=6712"
Replace-Text "847×8=6776" "374×4=1496"
Replace-Text "901×7=6307" "854×4=3416"
Replace-Text "775×9=6975" "607×2=1214"
Replace-Text "462×7=3234" "267×3=801"
Replace-Text "588×3=1764" "827×8=6616"
Replace-Text "238×6=1428" "882×2=1764"
Replace-Text "514×3=1542" "600×3=1800"
Replace-Text "745×5=3725" "551×2=1102"
Replace-Text "650×8=5200" "638×7=4466"
Replace-Text "816×8=6528" "561×8=4488"
Replace-Text "812×2=1624" "999×8=7992"
Replace-Text "609×2=1218" "925×9=8325"
Replace-Text "259×2=518" "243×5=1215"
Replace-Text "746×2=1492" "613×8=4904"
Replace-Text "129×3=387" "395×3=1185"
Replace-Text "872×7=6104" "307×3=921"
Replace-Text "954×4=3816" "392×4=1568"
Replace-Text "921×6=5526" "663×2=1326"
Replace-Text "803×8=6424" "631×8=5048"
